$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 1943.3334
$ws.Cells.Item(12, 9).Value = 2910
$ws.Cells.Item(12, 10).Value = 10
$ws.Cells.Item(12, 11).Value = 2910
$ws.Cells.Item(12, 12).Value = 10
$ws.Cells.Item(12, 13).Value = -2740
$ws.Cells.Item(12, 14).Value = -350

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 4123.2856
$ws.Cells.Item(106, 9).Value = 3858.08
$ws.Cells.Item(106, 11).Value = 3858.08
$ws.Cells.Item(106, 13).Value = -3227.08

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1334.5652
$ws.Cells.Item(112, 9).Value = 1333.3334
$ws.Cells.Item(112, 10).Value = 1334.75
$ws.Cells.Item(112, 11).Value = 4000.0002
$ws.Cells.Item(112, 12).Value = 4004.25
$ws.Cells.Item(112, 13).Value = -2892.0002
$ws.Cells.Item(112, 14).Value = -6220.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 7849
$ws.Cells.Item(116, 9).Value = 5759.6
$ws.Cells.Item(116, 11).Value = 5759.6
$ws.Cells.Item(116, 13).Value = -2317.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 1845.7273

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value = 60000
$ws.Cells.Item(133, 10).Value = 60000
$ws.Cells.Item(133, 12).Value = 60000
$ws.Cells.Item(133, 14).Value = -70120

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(136, 8).Value = 40000
$ws.Cells.Item(136, 10).Value = 40000
$ws.Cells.Item(136, 12).Value = 40000
$ws.Cells.Item(136, 14).Value = -50200

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3115.1707
$ws.Cells.Item(138, 9).Value = 1495.3636
$ws.Cells.Item(138, 10).Value = 3709.1
$ws.Cells.Item(138, 11).Value = 4486.0908
$ws.Cells.Item(138, 12).Value = 11127.3
$ws.Cells.Item(138, 13).Value = 653.9092000000001
$ws.Cells.Item(138, 14).Value = -21407.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 75215
$ws.Cells.Item(36, 9).Value = 85250.836
$ws.Cells.Item(36, 11).Value = 85250.836
$ws.Cells.Item(36, 13).Value = -84904.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 8553556
$ws.Cells.Item(94, 9).Value = 5652.75
$ws.Cells.Item(94, 11).Value = 5652.75
$ws.Cells.Item(94, 13).Value = -5201.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(45, 8).Value = 6500
$ws.Cells.Item(45, 9).Value = 6500
$ws.Cells.Item(45, 11).Value = 6500
$ws.Cells.Item(45, 13).Value = -5907

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(48, 8).Value = 26247.5
$ws.Cells.Item(48, 10).Value = 26247.5
$ws.Cells.Item(48, 12).Value = 26247.5
$ws.Cells.Item(48, 14).Value = -27199.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(56, 8).Value = 4000
$ws.Cells.Item(56, 9).Value = 4000
$ws.Cells.Item(56, 11).Value = 4000
$ws.Cells.Item(56, 13).Value = -3155

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(64, 8).Value = 50000
$ws.Cells.Item(64, 10).Value = 50000
$ws.Cells.Item(64, 12).Value = 50000
$ws.Cells.Item(64, 14).Value = -50496

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(67, 8).Value = 50000
$ws.Cells.Item(67, 10).Value = 50000
$ws.Cells.Item(67, 12).Value = 50000
$ws.Cells.Item(67, 14).Value = -51716

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(69, 8).Value = 20000
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 20000
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 20000
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(69, 14).Value = -21498

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(72, 8).Value = 20000
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 20000
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 60000
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(72, 14).Value = -67488

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(75, 8).Value = 40260
$ws.Cells.Item(75, 10).Value = 40260
$ws.Cells.Item(75, 12).Value = 40260
$ws.Cells.Item(75, 14).Value = -42256

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(78, 8).Value = 40260
$ws.Cells.Item(78, 10).Value = 40260
$ws.Cells.Item(78, 12).Value = 120780
$ws.Cells.Item(78, 14).Value = -130764

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(80, 8).Value = 28666.666
$ws.Cells.Item(80, 9).Value = 26000
$ws.Cells.Item(80, 11).Value = 26000
$ws.Cells.Item(80, 13).Value = -24877

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(82, 8).Value = 25181
$ws.Cells.Item(82, 10).Value = 25181
$ws.Cells.Item(82, 12).Value = 25181
$ws.Cells.Item(82, 14).Value = -25903

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(83, 8).Value = 28666.666
$ws.Cells.Item(83, 9).Value = 26000
$ws.Cells.Item(83, 11).Value = 78000
$ws.Cells.Item(83, 13).Value = -72384

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(85, 8).Value = 25181
$ws.Cells.Item(85, 10).Value = 25181
$ws.Cells.Item(85, 12).Value = 25181
$ws.Cells.Item(85, 14).Value = -27677

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(88, 8).Value = 8671.5
$ws.Cells.Item(88, 10).Value = 8671.5
$ws.Cells.Item(88, 12).Value = 8671.5
$ws.Cells.Item(88, 14).Value = -9483.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(91, 8).Value = 8671.5
$ws.Cells.Item(91, 10).Value = 8671.5
$ws.Cells.Item(91, 12).Value = 8671.5
$ws.Cells.Item(91, 14).Value = -11479.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(96, 8).Value = 8127.3076
$ws.Cells.Item(96, 10).Value = 8127.3076
$ws.Cells.Item(96, 12).Value = 8127.3076
$ws.Cells.Item(96, 14).Value = -13619.3076

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 1769.7
$ws.Cells.Item(105, 9).Value = 1769.7
$ws.Cells.Item(105, 11).Value = 1769.7
$ws.Cells.Item(105, 13).Value = -22.70000000000005

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 5438.636
$ws.Cells.Item(107, 9).Value = 709.5625
$ws.Cells.Item(107, 11).Value = 709.5625
$ws.Cells.Item(107, 13).Value = 1210.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2501.4211
$ws.Cells.Item(132, 10).Value = 3197.3333
$ws.Cells.Item(132, 12).Value = 9591.999899999999
$ws.Cells.Item(132, 14).Value = -14651.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 2012.4688
$ws.Cells.Item(107, 9).Value = 2340.125
$ws.Cells.Item(107, 10).Value = 1903.25
$ws.Cells.Item(107, 11).Value = 7020.375
$ws.Cells.Item(107, 12).Value = 5709.75
$ws.Cells.Item(107, 13).Value = -5100.375
$ws.Cells.Item(107, 14).Value = -9549.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 24204.5
$ws.Cells.Item(33, 10).Value = 24204.5
$ws.Cells.Item(33, 12).Value = 24204.5
$ws.Cells.Item(33, 14).Value = -24708.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(38, 8).Value = 21493.5
$ws.Cells.Item(38, 10).Value = 21493.5
$ws.Cells.Item(38, 12).Value = 21493.5
$ws.Cells.Item(38, 14).Value = -22419.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 20062.8
$ws.Cells.Item(40, 10).Value = 24999
$ws.Cells.Item(40, 12).Value = 24999
$ws.Cells.Item(40, 14).Value = -25301

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(55, 8).Value = 11183
$ws.Cells.Item(55, 10).Value = 14733
$ws.Cells.Item(55, 12).Value = 14733
$ws.Cells.Item(55, 14).Value = -15387

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64, 8).Value = 49500
$ws.Cells.Item(64, 10).Value = 49500
$ws.Cells.Item(64, 12).Value = 49500
$ws.Cells.Item(64, 14).Value = -49996

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(67, 8).Value = 49500
$ws.Cells.Item(67, 10).Value = 49500
$ws.Cells.Item(67, 12).Value = 49500
$ws.Cells.Item(67, 14).Value = -51216

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 20000
$ws.Cells.Item(126, 9).Value = 20000
$ws.Cells.Item(126, 11).Value = 60000
$ws.Cells.Item(126, 13).Value = -57530

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2625.4285
$ws.Cells.Item(46, 9).Value = 1500
$ws.Cells.Item(46, 11).Value = 1500
$ws.Cells.Item(46, 13).Value = -1312

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(127, 8).Value = 58799.4
$ws.Cells.Item(127, 10).Value = 58799.4
$ws.Cells.Item(127, 12).Value = 58799.4
$ws.Cells.Item(127, 14).Value = -68719.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3000
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 13).Value = -6530
